# Scheduled-runner update: refresh cached marketboard price/profit figures
# (currentAveragePrice / *NQ / *HQ / LevePriceNQ / LevePriceHQ / LeveProfitNQ
# / LeveProfitHQ columns H..N) across the per-job Leve tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 270.45456
$ws.Range("J55").Value = 636.5
$ws.Range("L55").Value = 636.5
$ws.Range("N55").Value = -1064.5
$ws.Range("H74").Value = 9360.154
$ws.Range("I74").Value = 6210.25
$ws.Range("K74").Value = 6210.25
$ws.Range("M74").Value = -5274.25
$ws.Range("H77").Value = 9360.154
$ws.Range("I77").Value = 6210.25
$ws.Range("K77").Value = 31051.25
$ws.Range("M77").Value = -26371.25
$ws.Range("H80").Value = 948
$ws.Range("I80").Value = 750
$ws.Range("J80").Value = 1032.8572
$ws.Range("K80").Value = 2250
$ws.Range("L80").Value = 3098.5716
$ws.Range("M80").Value = -1252
$ws.Range("N80").Value = -5094.571599999999
$ws.Range("H83").Value = 948
$ws.Range("I83").Value = 750
$ws.Range("J83").Value = 1032.8572
$ws.Range("K83").Value = 6750
$ws.Range("L83").Value = 9295.7148
$ws.Range("M83").Value = -1758
$ws.Range("N83").Value = -19279.7148
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120
$ws.Range("H137").Value = 1818.091
$ws.Range("I137").Value = 1561
$ws.Range("J137").Value = 2975
$ws.Range("K137").Value = 4683
$ws.Range("L137").Value = 8925
$ws.Range("M137").Value = -2133
$ws.Range("N137").Value = -14025
$ws.Range("H138").Value = 2584.0278
$ws.Range("I138").Value = 2795.5454
$ws.Range("J138").Value = 2490.96
$ws.Range("K138").Value = 8386.636200000001
$ws.Range("L138").Value = 7472.88
$ws.Range("M138").Value = -3246.636200000001
$ws.Range("N138").Value = -17752.88
$ws.Range("H141").Value = 2715.2222
$ws.Range("I141").Value = 2715.2222
$ws.Range("K141").Value = 8145.6666
$ws.Range("M141").Value = -2965.6666
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2850.4067
$ws.Range("I32").Value = 2727.138
$ws.Range("J32").Value = 10000
$ws.Range("K32").Value = 2727.138
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = -2440.138
$ws.Range("N32").Value = -10574
$ws.Range("H45").Value = 205199
$ws.Range("I45").Value = 669663.3
$ws.Range("K45").Value = 669663.3
$ws.Range("M45").Value = -669286.3
$ws.Range("H61").Value = 8401.258
$ws.Range("I61").Value = 7209.231
$ws.Range("K61").Value = 7209.231
$ws.Range("M61").Value = -6997.231
$ws.Range("H74").Value = 5162.227
$ws.Range("I74").Value = 4021.7646
$ws.Range("K74").Value = 4021.7646
$ws.Range("M74").Value = -3147.7646
$ws.Range("H77").Value = 5162.227
$ws.Range("I77").Value = 4021.7646
$ws.Range("K77").Value = 20108.823
$ws.Range("M77").Value = -15740.823
$ws.Range("H136").Value = 8401.258
$ws.Range("I136").Value = 7209.231
$ws.Range("K136").Value = 21627.693
$ws.Range("M136").Value = -19077.693
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2226.0625
$ws.Range("J94").Value = 4800
$ws.Range("L94").Value = 4800
$ws.Range("N94").Value = -5702
$ws.Range("H99").Value = 5690
$ws.Range("J99").Value = 6668.75
$ws.Range("L99").Value = 6668.75
$ws.Range("N99").Value = -9664.75
$ws.Range("H107").Value = 5250
$ws.Range("I107").Value = 5250
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 5250
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -3330
$ws.Range("N107").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4043.682
$ws.Range("I31").Value = 3259.3635
$ws.Range("J31").Value = 4828
$ws.Range("K31").Value = 3259.3635
$ws.Range("L31").Value = 4828
$ws.Range("M31").Value = -2964.3635
$ws.Range("N31").Value = -5418
$ws.Range("H34").Value = 4043.682
$ws.Range("I34").Value = 3259.3635
$ws.Range("J34").Value = 4828
$ws.Range("K34").Value = 3259.3635
$ws.Range("L34").Value = 4828
$ws.Range("M34").Value = -3057.3635
$ws.Range("N34").Value = -5232
$ws.Range("H99").Value = 4274.8667
$ws.Range("I99").Value = 4939.6665
$ws.Range("K99").Value = 4939.6665
$ws.Range("M99").Value = -3441.6665
$ws.Range("H122").Value = 6428.5713
$ws.Range("I122").Value = 5416.6665
$ws.Range("J122").Value = 12500
$ws.Range("K122").Value = 16249.9995
$ws.Range("L122").Value = 37500
$ws.Range("M122").Value = -13799.9995
$ws.Range("N122").Value = -42400
$ws.Range("H123").Value = 49666.668
$ws.Range("J123").Value = 49666.668
$ws.Range("L123").Value = 49666.668
$ws.Range("N123").Value = -59466.668
$ws.Range("H125").Value = 88108.664
$ws.Range("J125").Value = 88108.664
$ws.Range("L125").Value = 88108.664
$ws.Range("N125").Value = -93028.664
$ws.Range("H126").Value = 4274.8667
$ws.Range("I126").Value = 4939.6665
$ws.Range("K126").Value = 14818.9995
$ws.Range("M126").Value = -12348.9995
$ws.Range("H132").Value = 3695.125
$ws.Range("I132").Value = 3513.1428
$ws.Range("K132").Value = 10539.4284
$ws.Range("M132").Value = -8009.428400000001
$ws.Range("H134").Value = 4518.0386
$ws.Range("I134").Value = 3644.5417
$ws.Range("K134").Value = 10933.6251
$ws.Range("M134").Value = -8398.625100000001
$ws.Range("H138").Value = 95488.88
$ws.Range("J138").Value = 141090.8
$ws.Range("L138").Value = 141090.8
$ws.Range("N138").Value = -151370.8
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 140630.1
$ws.Range("J37").Value = 140630.1
$ws.Range("L37").Value = 421890.3
$ws.Range("N37").Value = -422114.3
$ws.Range("H128").Value = 1051068.5
$ws.Range("I128").Value = 1051068.5
$ws.Range("K128").Value = 3153205.5
$ws.Range("M128").Value = -3148225.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1335466.6
$ws.Range("I113").Value = 1335466.6
$ws.Range("K113").Value = 1335466.6
$ws.Range("M113").Value = -1333296.6
$ws.Range("H136").Value = 49289.617
$ws.Range("J136").Value = 49289.617
$ws.Range("L136").Value = 147868.851
$ws.Range("N136").Value = -152968.851
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3272.386
$ws.Range("I136").Value = 2854.3264
$ws.Range("J136").Value = 5833
$ws.Range("K136").Value = 8562.9792
$ws.Range("L136").Value = 17499
$ws.Range("M136").Value = -6012.9792
$ws.Range("N136").Value = -22599
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 21000
$ws.Range("I44").Value = 15000
$ws.Range("K44").Value = 15000
$ws.Range("M44").Value = -14446
$ws.Range("H46").Value = 106714.5
$ws.Range("J46").Value = 106714.5
$ws.Range("L46").Value = 106714.5
$ws.Range("N46").Value = -107176.5
$ws.Range("H47").Value = 20000
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H119").Value = 75000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 75000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 75000
$ws.Range("M119").ClearContents()
$ws.Range("N119").Value = -84676
$ws.Range("H126").Value = 2960.1924
$ws.Range("I126").Value = 2866.2778
$ws.Range("K126").Value = 8598.8334
$ws.Range("M126").Value = -6128.8334
$ws.Range("H132").Value = 3452.8
$ws.Range("I132").Value = 3721.5312
$ws.Range("K132").Value = 11164.5936
$ws.Range("M132").Value = -8634.5936
$ws.Range("H134").Value = 106714.5
$ws.Range("J134").Value = 106714.5
$ws.Range("L134").Value = 320143.5
$ws.Range("N134").Value = -325213.5
$ws.Range("H136").Value = 9700.454
$ws.Range("I136").Value = 6562.1665
$ws.Range("K136").Value = 19686.4995
$ws.Range("M136").Value = -17136.4995
